$wb = $excel.ActiveWorkbook

# The same "想去人数" (interested-count) updates apply to both the
# "展览" and "全部类型" worksheets, which mirror the same event rows.
$sheetNames = @("展览", "全部类型")

# Map of cell address -> new value
$updates = @{
    "F3"  = 1713
    "F7"  = 11967
    "F10" = 476
    "F13" = 855
    "F15" = 13449
    "F23" = 239
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
